$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 15873550
$ws.Range("I33").Value = 23810092
$ws.Range("K33").Value = 23810092
$ws.Range("M33").Value = -23809863
$ws.Range("H80").Value = 1016.0625
$ws.Range("I80").Value = 231.57143
$ws.Range("J80").Value = 1626.2222
$ws.Range("K80").Value = 694.71429
$ws.Range("L80").Value = 4878.6666
$ws.Range("M80").Value = 303.28571
$ws.Range("N80").Value = -6874.6666
$ws.Range("H83").Value = 1016.0625
$ws.Range("I83").Value = 231.57143
$ws.Range("J83").Value = 1626.2222
$ws.Range("K83").Value = 2084.14287
$ws.Range("L83").Value = 14635.9998
$ws.Range("M83").Value = 2907.85713
$ws.Range("N83").Value = -24619.9998
$ws.Range("H98").Value = 5732.773
$ws.Range("I98").Value = 5986.7144
$ws.Range("J98").Value = 400
$ws.Range("K98").Value = 5986.7144
$ws.Range("L98").Value = 400
$ws.Range("M98").Value = -4488.7144
$ws.Range("N98").Value = -3396
$ws.Range("H106").Value = 16701957
$ws.Range("I106").Value = 38840.555
$ws.Range("K106").Value = 38840.555
$ws.Range("M106").Value = -38209.555
$ws.Range("H116").Value = 7333.3335
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 10000
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 10000
$ws.Range("M116").Value = 1442
$ws.Range("N116").Value = -16884
$ws.Range("H122").Value = 5732.773
$ws.Range("I122").Value = 5986.7144
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 17960.1432
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = -15510.1432
$ws.Range("N122").Value = -6100
$ws.Range("H137").Value = 50002216
$ws.Range("I137").Value = 1488.6
$ws.Range("J137").Value = 200004400
$ws.Range("K137").Value = 4465.799999999999
$ws.Range("L137").Value = 600013200
$ws.Range("M137").Value = -1915.799999999999
$ws.Range("N137").Value = -600018300
$ws.Range("H138").Value = 3046.4268
$ws.Range("I138").Value = 2844.2778
$ws.Range("J138").Value = 3103.2812
$ws.Range("K138").Value = 8532.8334
$ws.Range("L138").Value = 9309.8436
$ws.Range("M138").Value = -3392.8334
$ws.Range("N138").Value = -19589.8436
$ws.Range("H141").Value = 2108.3333
$ws.Range("I141").Value = 1435.8334
$ws.Range("J141").Value = 3261.1904
$ws.Range("K141").Value = 4307.5002
$ws.Range("L141").Value = 9783.5712
$ws.Range("M141").Value = 872.4997999999996
$ws.Range("N141").Value = -20143.5712
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17309.041
$ws.Range("I32").Value = 17068.28
$ws.Range("J32").Value = 19034.5
$ws.Range("K32").Value = 17068.28
$ws.Range("L32").Value = 19034.5
$ws.Range("M32").Value = -16781.28
$ws.Range("N32").Value = -19608.5
$ws.Range("H61").Value = 2166.125
$ws.Range("I61").Value = 2124.889
$ws.Range("J61").Value = 2190.8667
$ws.Range("K61").Value = 2124.889
$ws.Range("L61").Value = 2190.8667
$ws.Range("M61").Value = -1912.889
$ws.Range("N61").Value = -2614.8667
$ws.Range("H132").Value = 6341.08
$ws.Range("I132").Value = 7340.028
$ws.Range("J132").Value = 3772.3572
$ws.Range("K132").Value = 22020.084
$ws.Range("L132").Value = 11317.0716
$ws.Range("M132").Value = -19490.084
$ws.Range("N132").Value = -16377.0716
$ws.Range("H136").Value = 2166.125
$ws.Range("I136").Value = 2124.889
$ws.Range("J136").Value = 2190.8667
$ws.Range("K136").Value = 6374.667
$ws.Range("L136").Value = 6572.6001
$ws.Range("M136").Value = -3824.667
$ws.Range("N136").Value = -11672.6001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3407.8333
$ws.Range("I31").Value = 2746
$ws.Range("J31").Value = 5922.8
$ws.Range("K31").Value = 2746
$ws.Range("L31").Value = 5922.8
$ws.Range("M31").Value = -2451
$ws.Range("N31").Value = -6512.8
$ws.Range("H34").Value = 3407.8333
$ws.Range("I34").Value = 2746
$ws.Range("J34").Value = 5922.8
$ws.Range("K34").Value = 2746
$ws.Range("L34").Value = 5922.8
$ws.Range("M34").Value = -2544
$ws.Range("N34").Value = -6326.8
$ws.Range("H58").Value = 1484.1177
$ws.Range("I58").Value = 1498.3636
$ws.Range("J58").Value = 1014
$ws.Range("K58").Value = 1498.3636
$ws.Range("L58").Value = 1014
$ws.Range("M58").Value = -1295.3636
$ws.Range("N58").Value = -1420
$ws.Range("H68").Value = 37949.832
$ws.Range("J68").Value = 39539.8
$ws.Range("L68").Value = 39539.8
$ws.Range("N68").Value = -41037.8
$ws.Range("H71").Value = 37949.832
$ws.Range("J71").Value = 39539.8
$ws.Range("L71").Value = 118619.4
$ws.Range("N71").Value = -126107.4
$ws.Range("H99").Value = 2352.5264
$ws.Range("I99").Value = 2483
$ws.Range("J99").Value = 2292.3076
$ws.Range("K99").Value = 2483
$ws.Range("L99").Value = 2292.3076
$ws.Range("M99").Value = -985
$ws.Range("N99").Value = -5288.3076
$ws.Range("H126").Value = 2352.5264
$ws.Range("I126").Value = 2483
$ws.Range("J126").Value = 2292.3076
$ws.Range("K126").Value = 7449
$ws.Range("L126").Value = 6876.9228
$ws.Range("M126").Value = -4979
$ws.Range("N126").Value = -11816.9228
$ws.Range("H136").Value = 1484.1177
$ws.Range("I136").Value = 1498.3636
$ws.Range("J136").Value = 1014
$ws.Range("K136").Value = 4495.0908
$ws.Range("L136").Value = 3042
$ws.Range("M136").Value = -1945.0908
$ws.Range("N136").Value = -8142
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2181350.8
$ws.Range("I131").Value = 12894.875
$ws.Range("J131").Value = 2584784.2
$ws.Range("K131").Value = 38684.625
$ws.Range("L131").Value = 7754352.600000001
$ws.Range("M131").Value = -33644.625
$ws.Range("N131").Value = -7764432.600000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4758.5713
$ws.Range("I132").Value = 4985.5483
$ws.Range("J132").Value = 2999.5
$ws.Range("K132").Value = 14956.6449
$ws.Range("L132").Value = 8998.5
$ws.Range("M132").Value = -12426.6449
$ws.Range("N132").Value = -14058.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7429.9395
$ws.Range("I132").Value = 11858.353
$ws.Range("J132").Value = 2724.75
$ws.Range("K132").Value = 35575.05899999999
$ws.Range("L132").Value = 8174.25
$ws.Range("M132").Value = -33045.05899999999
$ws.Range("N132").Value = -13234.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1794.5807
$ws.Range("I132").Value = 1611.9286
$ws.Range("J132").Value = 3499.3333
$ws.Range("K132").Value = 4835.7858
$ws.Range("L132").Value = 10497.9999
$ws.Range("M132").Value = -2305.7858
$ws.Range("N132").Value = -15557.9999
$ws.Range("H141").Value = 70000
$ws.Range("J141").Value = 70000
$ws.Range("L141").Value = 70000
$ws.Range("N141").Value = -80360
